$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeGetResults")

$ws.Cells.Item(3, 1).Value = "https://www.youtube.com/watch?v=j1OrjcZyrhg&list=PLQOaTSbfxUtCrKs0nicOg2npJQYSPGO9r&index=26"
$ws.Cells.Item(3, 2).Value = 200

$ws.Cells.Item(4, 1).Value = "https://github.com/Okarpets/Sitest_applicatio"
$ws.Cells.Item(4, 2).Value = 404

$ws.Cells.Item(5, 1).Value = "https://github.com/Okarpets/Sitest_application"
$ws.Cells.Item(5, 2).Value = 200
